# Added listener and calculations for ECU CAN stream.
# Updates the "Typography" and "Translation" sheets of the TouchGFX
# texts workbook:
#   - Typography: change the "Default" typography to a bold 18px font,
#     and add two new typographies ("rpm_letters" / "rpm_descriptor")
#     used to render the RPM gauge text.
#   - Translation: rename a handful of labels (ECT -> EGT, OILT ->
#     OIL. T, OILP -> OIL. P), update the demo/preview value, repoint a
#     couple of rows at the Default typography, and add three new
#     translation rows for the RPM widgets.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("Typography")
$ws2 = $wb.Worksheets.Item("Translation")

# xlPasteFormats - used below to strip the "Text" number-format style
# that Excel applies when we force a numeric-looking literal (e.g.
# "740"/"9500") to be stored as text, so the cell keeps inheriting the
# plain column style instead of picking up a one-off style record.
$xlPasteFormats = -4122

# --- Typography sheet ---------------------------------------------------

# Row 4 ("Default" typography): switch from regular to bold, 20 -> 18 px.
$ws1.Cells.Item(4, 3).Value = "Asap-Bold.ttf"
$ws1.Cells.Item(4, 4).Value = 18

# Row 10: new "rpm_letters" typography (bold, 21px).
$ws1.Cells.Item(10, 2).Value = "rpm_letters"
$ws1.Cells.Item(10, 3).Value = "Asap-Bold.ttf"
$ws1.Cells.Item(10, 4).Value = 21
$ws1.Cells.Item(10, 5).Value = 4
$ws1.Cells.Item(10, 6).Value = "?"
$ws1.Cells.Item(10, 7).Value = "0-9"
$ws1.Cells.Item(10, 9).Value = "0-9"

# Row 11: new "rpm_descriptor" typography (bold, 10px).
$ws1.Cells.Item(11, 2).Value = "rpm_descriptor"
$ws1.Cells.Item(11, 3).Value = "Asap-Bold.ttf"
$ws1.Cells.Item(11, 4).Value = 10
$ws1.Cells.Item(11, 5).Value = 4
$ws1.Cells.Item(11, 6).Value = "?"

# --- Translation sheet ---------------------------------------------------

# Relabel / update existing rows. E28 and further numeric-looking
# labels must stay text, so force the format before writing and then
# restore the plain (un-styled) format from an untouched sibling cell.
$ws2.Cells.Item(28, 5).NumberFormat = "@"
$ws2.Cells.Item(28, 5).Value = "740"
$ws2.Cells.Item(29, 5).Copy()
$ws2.Cells.Item(28, 5).PasteSpecial($xlPasteFormats)

$ws2.Cells.Item(30, 5).Value = "EGT"
$ws2.Cells.Item(33, 5).Value = "OIL. T"
$ws2.Cells.Item(34, 5).Value = "OIL. P"

# Repoint these rows at the "Default" typography instead of "small".
$ws2.Cells.Item(42, 3).Value = "Default"
$ws2.Cells.Item(43, 3).Value = "Default"

# New rows for the RPM gauge widgets.
$ws2.Cells.Item(44, 2).Value = "SingleUseId46"
$ws2.Cells.Item(44, 3).Value = "rpm_letters"
$ws2.Cells.Item(44, 4).Value = "Right"
$ws2.Cells.Item(44, 5).Value = "<value> "
$ws2.Cells.Item(44, 6).Value = "LTR"

$ws2.Cells.Item(45, 2).Value = "SingleUseId47"
$ws2.Cells.Item(45, 3).Value = "rpm_letters"
$ws2.Cells.Item(45, 4).Value = "Left"
$ws2.Cells.Item(45, 5).NumberFormat = "@"
$ws2.Cells.Item(45, 5).Value = "9500"
$ws2.Cells.Item(44, 5).Copy()
$ws2.Cells.Item(45, 5).PasteSpecial($xlPasteFormats)
$ws2.Cells.Item(45, 6).Value = "LTR"

$ws2.Cells.Item(46, 2).Value = "SingleUseId50"
$ws2.Cells.Item(46, 3).Value = "rpm_descriptor"
$ws2.Cells.Item(46, 4).Value = "Left"
$ws2.Cells.Item(46, 5).Value = "RPM"
$ws2.Cells.Item(46, 6).Value = "LTR"

$excel.CutCopyMode = $false
